$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Sdate values for rows 2 and 4 (same text value as before -> new dates)
$ws.Range("C2").Value = "2023-11-05T00:00:00"
$ws.Range("C4").Value = "2023-11-03T00:00:00"

# Update the active selection to E17
$ws.Range("E17").Select()
